$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values in the source are plain text that sometimes happens to look
# numeric (e.g. "243.50", "1.0000"). Prefix those with a literal leading quote,
# the same trick the Excel UI uses to force "store as text", so the value
# round-trips exactly instead of being auto-coerced into a number and losing
# formatting like trailing zeros. Values with two dots (e.g. "30.293.12") are
# never valid numbers so they do not need this treatment.
$ws.Range("D2").Value = "30.293.12"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "1.866.60"
$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'243.50"
$ws.Range("E5").Value = "  -2.40%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.4727"
$ws.Range("E8").Value = "  -2.51%  "
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("D11").Value = "'0.07796"
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").Value = "'97.40"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "1.867.30"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "'0.7213"
$ws.Range("E14").Value = "  -2.62%  "
$ws.Range("D15").Value = "'5.147"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").Value = "'280.00"
$ws.Range("E16").Value = "  +1.72%  "
$ws.Range("D17").Value = "30.278.49"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").Value = "'13.01"
$ws.Range("E18").Value = "  -1.38%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").Value = "2.111.03"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("D23").Value = "'5.232"
$ws.Range("E23").Value = "  -2.20%  "
$ws.Range("D24").Value = "'6.261"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").Value = "'162.48"
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("D26").Value = "'9.000"
$ws.Range("E26").Value = "  -2.69%  "
$ws.Range("D27").Value = "'18.69"
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("E28").Value = "  -2.15%  "
$ws.Range("D29").Value = "'0.09619"
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("E30").Value = "  -2.24%  "
$ws.Range("D31").Value = "'1.478"
$ws.Range("E31").Value = "  -1.78%  "
$ws.Range("D32").Value = "'4.216"
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("D34").Value = "'0.04783"
$ws.Range("E34").Value = "  -2.43%  "
$ws.Range("D35").Value = "'1.117"
$ws.Range("E35").Value = "  -0.83%  "
$ws.Range("D36").Value = "'0.6846"
$ws.Range("E36").Value = "  -2.32%  "
$ws.Range("D37").Value = "'2.714"
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("D38").Value = "'0.01887"
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("D39").Value = "'2.838"
$ws.Range("E39").Value = "  +1.62%  "
$ws.Range("D40").Value = "'75.18"
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("D41").Value = "'6.203"
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("E42").Value = "  -4.52%  "
$ws.Range("D43").Value = "'0.4210"
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("D44").Value = "'0.9999"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'0.8250"
$ws.Range("E45").Value = "  -2.14%  "
$ws.Range("D46").Value = "'100.47"
$ws.Range("E46").Value = "  -2.24%  "
$ws.Range("D47").Value = "'9.651"
$ws.Range("E47").Value = "  +2.25%  "
$ws.Range("D48").Value = "'6.947"
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("D49").Value = "'34.91"
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("D50").Value = "'0.05766"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").Value = "'882.95"
$ws.Range("E51").Value = "  -3.85%  "
